$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22 / 23 swap: Litecoin <-> Dai (data reordered by the source ranking) ---
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "71.14"
$ws.Range("E23").Value = "  +3.35%  "

# --- Price / volume refresh for the remaining rows ---
$ws.Range("D2").Value = "68.529.82"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "2.519.82"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("D5").Value = "591.99"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "176.30"
$ws.Range("E6").Value = "  +1.85%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.517"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("D11").Value = "5.00"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").Value = "2.944.95"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "25.78"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "68.076.30"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").Value = "0.0000170"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "2.491.97"
$ws.Range("E17").Value = "  -1.84%  "
$ws.Range("D18").Value = "10.99"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").Value = "7.47"
$ws.Range("E19").Value = "  -1.68%  "
$ws.Range("D20").Value = "351.73"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "4.09"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("D24").Value = "4.23"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "1.73"
$ws.Range("E25").Value = "  -4.96%  "
$ws.Range("D26").Value = "9.03"
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("D27").Value = "2.587.99"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "0.0₃0901"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").Value = "506.13"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "7.84"
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("D36").Value = "162.40"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "18.36"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("E40").Value = "  +3.92%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").Value = "4.84"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("D44").Value = "2.41"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("D45").Value = "149.87"
$ws.Range("E45").Value = "  +4.89%  "
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "0.519"
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("D49").Value = "0.0738"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").Value = "0.580"
$ws.Range("E51").Value = "  -0.52%  "
